# Auto-generated: apply updated market-price figures scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1809.0698
$ws.Range("I132").Value = 1492.2821
$ws.Range("K132").Value = 4476.846299999999
$ws.Range("M132").Value = -1946.846299999999
$ws.Range("H135").Value = 84063.414
$ws.Range("I135").Value = 797.3333
$ws.Range("J135").Value = 111818.78
$ws.Range("K135").Value = 7175.9997
$ws.Range("L135").Value = 1006369.02
$ws.Range("M135").Value = -4640.9997
$ws.Range("N135").Value = -1011439.02

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23834.094
$ws.Range("I32").Value = 5061.6924
$ws.Range("K32").Value = 5061.6924
$ws.Range("M32").Value = -4774.6924
$ws.Range("H39").Value = 4739.75
$ws.Range("I39").Value = 4739.75
$ws.Range("K39").Value = 4739.75
$ws.Range("M39").Value = -4219.75
$ws.Range("H45").Value = 4680.45
$ws.Range("J45").Value = 3374.2144
$ws.Range("L45").Value = 3374.2144
$ws.Range("N45").Value = -4128.2144
$ws.Range("H61").Value = 1900.8572
$ws.Range("I61").Value = 1900.8572
$ws.Range("K61").Value = 1900.8572
$ws.Range("M61").Value = -1688.8572
$ws.Range("H132").Value = 1229.9412
$ws.Range("I132").Value = 1199.1277
$ws.Range("K132").Value = 3597.3831
$ws.Range("M132").Value = -1067.3831
$ws.Range("H136").Value = 1900.8572
$ws.Range("I136").Value = 1900.8572
$ws.Range("K136").Value = 5702.571599999999
$ws.Range("M136").Value = -3152.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13880.061
$ws.Range("I31").Value = 3777.4443
$ws.Range("J31").Value = 26003.2
$ws.Range("K31").Value = 3777.4443
$ws.Range("L31").Value = 26003.2
$ws.Range("M31").Value = -3482.4443
$ws.Range("N31").Value = -26593.2
$ws.Range("H34").Value = 13880.061
$ws.Range("I34").Value = 3777.4443
$ws.Range("J34").Value = 26003.2
$ws.Range("K34").Value = 3777.4443
$ws.Range("L34").Value = 26003.2
$ws.Range("M34").Value = -3575.4443
$ws.Range("N34").Value = -26407.2
$ws.Range("H35").Value = 2858.0908
$ws.Range("I35").Value = 1156.6666
$ws.Range("K35").Value = 1156.6666
$ws.Range("M35").Value = -862.6666
$ws.Range("H58").Value = 2173.158
$ws.Range("I58").Value = 2224.8
$ws.Range("J58").Value = 1979.5
$ws.Range("K58").Value = 2224.8
$ws.Range("L58").Value = 1979.5
$ws.Range("M58").Value = -2021.8
$ws.Range("N58").Value = -2385.5
$ws.Range("H136").Value = 2173.158
$ws.Range("I136").Value = 2224.8
$ws.Range("J136").Value = 1979.5
$ws.Range("K136").Value = 6674.400000000001
$ws.Range("L136").Value = 5938.5
$ws.Range("M136").Value = -4124.400000000001
$ws.Range("N136").Value = -11038.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 850
$ws.Range("H27").Value = 850
$ws.Range("H39").Value = 8568.4375
$ws.Range("J39").Value = 8568.4375
$ws.Range("L39").Value = 25705.3125
$ws.Range("N39").Value = -26293.3125
$ws.Range("H122").Value = 346.66666
$ws.Range("I122").Value = 285.57144
$ws.Range("J122").Value = 377.2143
$ws.Range("K122").Value = 2570.14296
$ws.Range("L122").Value = 3394.9287
$ws.Range("M122").Value = -120.1429600000001
$ws.Range("N122").Value = -8294.9287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 53499.25
$ws.Range("J47").Value = 53499.25
$ws.Range("L47").Value = 53499.25
$ws.Range("N47").Value = -54635.25
$ws.Range("H58").Value = 42999
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 42999
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 42999
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -43553
$ws.Range("H122").Value = 1464.8
$ws.Range("I122").Value = 1531
$ws.Range("K122").Value = 4593
$ws.Range("M122").Value = -2143
$ws.Range("H126").Value = 3463.6428
$ws.Range("I126").Value = 2856.2856
$ws.Range("K126").Value = 8568.856800000001
$ws.Range("M126").Value = -6098.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2783.3333
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 3200
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 3200
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -3790
$ws.Range("H27").Value = 2783.3333
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 3200
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 3200
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -3414
$ws.Range("H40").Value = 3690.2104
$ws.Range("I40").Value = 2942
$ws.Range("K40").Value = 2942
$ws.Range("M40").Value = -2806
$ws.Range("H43").Value = 11622594
$ws.Range("H82").Value = 1646.7646
$ws.Range("I82").Value = 1445.4546
$ws.Range("J82").Value = 2015.8334
$ws.Range("K82").Value = 1445.4546
$ws.Range("L82").Value = 2015.8334
$ws.Range("M82").Value = -1084.4546
$ws.Range("N82").Value = -2737.8334
$ws.Range("H85").Value = 1646.7646
$ws.Range("I85").Value = 1445.4546
$ws.Range("J85").Value = 2015.8334
$ws.Range("K85").Value = 1445.4546
$ws.Range("L85").Value = 2015.8334
$ws.Range("M85").Value = -197.4546
$ws.Range("N85").Value = -4511.8334
$ws.Range("H132").Value = 3011.5134
$ws.Range("I132").Value = 2279.6956
$ws.Range("K132").Value = 6839.0868
$ws.Range("M132").Value = -4309.0868
$ws.Range("H136").Value = 4984.1816
$ws.Range("I136").Value = 4243.6
$ws.Range("K136").Value = 12730.8
$ws.Range("M136").Value = -10180.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 559.2222
$ws.Range("I81").Value = 559.2222
$ws.Range("K81").Value = 1118.4444
$ws.Range("M81").Value = -57.44440000000009
$ws.Range("H84").Value = 559.2222
$ws.Range("I84").Value = 559.2222
$ws.Range("K84").Value = 5592.222000000001
$ws.Range("M84").Value = -288.2220000000007
$ws.Range("H126").Value = 2411.611
$ws.Range("I126").Value = 2121.3572
$ws.Range("K126").Value = 6364.071599999999
$ws.Range("M126").Value = -3894.071599999999
$ws.Range("H132").Value = 1279178.5
$ws.Range("I132").Value = 1468273.6
$ws.Range("K132").Value = 4404820.800000001
$ws.Range("M132").Value = -4402290.800000001
$ws.Range("H136").Value = 664.75
$ws.Range("I136").Value = 664.75
$ws.Range("K136").Value = 1994.25
$ws.Range("M136").Value = 555.75

